# "envoi de quelques documents"
# Fill in the work log for Iteration #1 (rows 14-21: date / task / hours),
# then leave the workbook scrolled to, and selected on, that sheet - as if
# the author had just finished typing the entries and saved.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Estimation")
$ws2 = $wb.Worksheets.Item("Iteration #1")

# --- Iteration #1 work log (A14:C21) -------------------------------------
# A14 already carries a date number format; A15:A21 don't yet, so give them
# one (this is what Excel itself does the first time you type a date into a
# previously-unformatted cell). Format A15 then fan the format out to the
# rest of the column so they all share one style, same as a manual fill.
$ws2.Range("A15").NumberFormat = "mm-dd-yy"
$ws2.Range("A15").Copy()
$ws2.Range("A16:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A14").Value = Get-Date -Year 2017 -Month 1  -Day 30 -Hour 0 -Minute 0 -Second 0
$ws2.Range("B14").Value = "commencer les  activity et le java"
$ws2.Range("C14").Value = 3

$ws2.Range("A15").Value = Get-Date -Year 2017 -Month 1  -Day 31 -Hour 0 -Minute 0 -Second 0
$ws2.Range("B15").Value = "mise en place du serveur apache, php, mysql et phpmyadmin"
$ws2.Range("C15").Value = 2

$ws2.Range("A16").Value = Get-Date -Year 2017 -Month 2  -Day 6  -Hour 0 -Minute 0 -Second 0
$ws2.Range("B16").Value = "création de la bd"
$ws2.Range("C16").Value = 3

$ws2.Range("A17").Value = Get-Date -Year 2017 -Month 2  -Day 7  -Hour 0 -Minute 0 -Second 0
$ws2.Range("B17").Value = "tentative de connexion à la bd (échec)"
$ws2.Range("C17").Value = 2

$ws2.Range("A18").Value = Get-Date -Year 2017 -Month 2  -Day 9  -Hour 0 -Minute 0 -Second 0
$ws2.Range("B18").Value = "documentation"
$ws2.Range("C18").Value = 3

$ws2.Range("A19").Value = Get-Date -Year 2017 -Month 2  -Day 13 -Hour 0 -Minute 0 -Second 0
$ws2.Range("B19").Value = "tentative de requête à l'aide d'android(échec)"
$ws2.Range("C19").Value = 3

$ws2.Range("A20").Value = Get-Date -Year 2017 -Month 2  -Day 14 -Hour 0 -Minute 0 -Second 0
$ws2.Range("B20").Value = "documentation concernant android et php"
$ws2.Range("C20").Value = 5

$ws2.Range("A21").Value = Get-Date -Year 2017 -Month 2  -Day 17 -Hour 0 -Minute 0 -Second 0
$ws2.Range("B21").Value = "documentation concernant android et php"
$ws2.Range("C21").Value = 2

# --- Window state: the author re-opened the book on the "Iteration #1" tab,
# scrolled down to the log, and left the cursor on B23 --------------------
$ws2.Activate()
$ws2.Range("B23").Select()
